$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the "Last compiled" date: 2024-08-26 -> 2024-09-04
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2024-08-26", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-04", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert a new "There will be no **final exam** for the class." paragraph
#    right after the StatsChats paragraph (the one ending in "...during the
#    semester.") and before the "Grading" Heading3 paragraph.
# ---------------------------------------------------------------------------
$statsChatsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "during the semester") {
        $statsChatsPara = $p
    }
}

$gradingPara = $statsChatsPara.Next()

# Grab the (single) character at the very start of the "Grading" paragraph so
# we operate on a non-collapsed Range - this keeps paragraph boundaries
# intact when injecting raw WordprocessingML via InsertXML.
$graduatingStart = $gradingPara.Range.Start
$placeholderRange = $d.Range($graduatingStart, $graduatingStart + 1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">There will be no</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:bCs/><w:b/></w:rPr><w:t xml:space="preserve">final exam</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">for the class.</w:t></w:r></w:p>'
$gradingParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t xml:space="preserve">Grading</w:t></w:r></w:p>'

$placeholderRange.InsertXML($newParaXml + $gradingParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 3. Change the "Download the Syllabus" heading from Heading2 to Heading3.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Download the Syllabus") {
        $p.Range.ParagraphFormat.Style = "Heading3"
    }
}
